# Update "想去人数" (F column) values for a few rows on the
# "展览" and "全部类型" sheets to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 22
$ws1.Range("F7").Value = 3405
$ws1.Range("F9").Value = 4083
$ws1.Range("F12").Value = 42

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 22
$ws4.Range("F8").Value = 3405
$ws4.Range("F10").Value = 4083
$ws4.Range("F13").Value = 42
